$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename existing shared string used by D12 from "Billboard" to "Buliding"
$ws.Range("D12").Value = "Buliding"

# Add new row 13 data
$ws.Range("B13").Value = 9
$ws.Range("D13").Value = "Enemy"
$ws.Range("E13").Value = 60
$ws.Range("F13").Value = 60
$ws.Range("G13").Value = 60
$ws.Range("H13").Value = 0

# Update selection to match diff (K9 selected)
$ws.Range("K9").Select()
